$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Timesheet" ---
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Row 2: Osborn / PTO -> Layne / Regular, rate+total zeroed
$ws1.Range("B2").Value = "Layne"
$ws1.Range("D2").Value = "Regular"
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = 0

# Row 3: Caputo (Insp.) -> Campbell, rate+total zeroed
$ws1.Range("B3").Value = "Campbell"
$ws1.Range("E3").Value = 0
$ws1.Range("F3").Value = 0

# Row 4: Funke (Insp.) -> Cummin, rate+total zeroed
$ws1.Range("B4").Value = "Cummin"
$ws1.Range("E4").Value = 0
$ws1.Range("F4").Value = 0

# Row 5: Winn -> McClure, rate+total zeroed
$ws1.Range("B5").Value = "McClure"
$ws1.Range("E5").Value = 0
$ws1.Range("F5").Value = 0

# Row 6: Muncey (Maint. Items) -> Hunter, rate+total zeroed
$ws1.Range("B6").Value = "Hunter"
$ws1.Range("E6").Value = 0
$ws1.Range("F6").Value = 0

# Subtotal rows: totals zeroed out (reverting seeded sample rate/total)
$ws1.Range("F8").Value = 0
$ws1.Range("F11").Value = 0
$ws1.Range("F13").Value = 0

# --- Sheet 2: "Jason Schema" ---
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Employee ID changed for all data rows
$ws2.Range("B2").Value = "emp_c2dcy26q"
$ws2.Range("B3").Value = "emp_c2dcy26q"
$ws2.Range("B4").Value = "emp_c2dcy26q"
$ws2.Range("B5").Value = "emp_c2dcy26q"
$ws2.Range("B6").Value = "emp_c2dcy26q"

# Row 2: Osborn/PTO -> Layne/Regular, rate+total zeroed, note cleared
$ws2.Range("D2").Value = "Layne"
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0
$ws2.Range("H2").Value = "Regular"
$ws2.Range("I2").Value = ""

# Row 3: Caputo (Insp.) -> Campbell, rate+total zeroed, note cleared
$ws2.Range("D3").Value = "Campbell"
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 0
$ws2.Range("I3").Value = ""

# Row 4: Funke (Insp.) -> Cummin, rate+total zeroed, note cleared
$ws2.Range("D4").Value = "Cummin"
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 0
$ws2.Range("I4").Value = ""

# Row 5: Winn -> McClure, rate+total zeroed, note cleared
$ws2.Range("D5").Value = "McClure"
$ws2.Range("F5").Value = 0
$ws2.Range("G5").Value = 0
$ws2.Range("I5").Value = ""

# Row 6: Muncey (Maint. Items) -> Hunter, rate+total zeroed, note cleared
$ws2.Range("D6").Value = "Hunter"
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 0
$ws2.Range("I6").Value = ""
